# Insert a new row above the old row 2 ("Total time taken for the ride")
# for the new "Date and Time" metric, and a new row above the old row 35
# ("Idling time percentage") for the new "Cycle Count of battery" metric.
# All the rows in between (and after) shift down accordingly, which matches
# the diff's renumbering of rows 2-44 -> 3-46 (and 35-44 -> 37-46 after the
# second insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for "Date and Time" right after the header row.
$ws.Rows.Item(2).Insert()

# 2) Make room for "Cycle Count of battery" right after
#    "Electricity consumption units(kW)" (old row 34, now row 35) and right
#    before "Idling time percentage" (old row 35, now row 36).
$ws.Rows.Item(36).Insert()

# Populate the new "Date and Time" row.
$ws.Range("A2").Value = "Date and Time"
$ws.Range("B2").Value = "2024-03-11 13:17:50.432000 to 2024-03-11 15:29:55.824000"
$ws.Range("C2").Value = "2024-03-11 17:42:37.788000 to 2024-03-11 20:01:41.810000"
$ws.Range("D2").Value = "2024-03-11 16:37:01.599000 to 2024-03-11 17:36:23.668000"
$ws.Range("E2").Value = "2024-03-11 15:31:19.824000 to 2024-03-11 16:36:21.599000"

# Populate the new "Cycle Count of battery" row.
$ws.Range("A36").Value = "Cycle Count of battery"
$ws.Range("B36").Value = 75
$ws.Range("C36").Value = 136
$ws.Range("D36").Value = 42
$ws.Range("E36").Value = 114
